$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.697.74'
$ws.Range("E2").Value = '  -0.90%  '
$ws.Range("D3").Value = '1.629.60'
$ws.Range("E3").Value = '  -1.05%  '
$ws.Range("E4").Value = '  -0.34%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.21'
$ws.Range("E5").Value = '  -1.01%  '
$ws.Range("E6").Value = '  -1.02%  '
$ws.Range("E7").Value = '  -0.24%  '
$ws.Range("E8").Value = '  -0.70%  '
$ws.Range("E9").Value = '  -1.48%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.47'
$ws.Range("E10").Value = '  -4.98%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0782'
$ws.Range("E11").Value = '  +0.00%  '
$ws.Range("E12").Value = '  -0.99%  '
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '1.854.41'
$ws.Range("E13").Value = '  -1.05%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.619.04'
$ws.Range("E14").Value = '  -1.69%  '
$ws.Range("E15").Value = '  -1.89%  '
$ws.Range("D16").Value = '0.0₃0762'
$ws.Range("E16").Value = '  -1.18%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.12'
$ws.Range("E17").Value = '  -0.55%  '
$ws.Range("D18").Value = '25.722.22'
$ws.Range("E18").Value = '  -0.89%  '
$ws.Range("E19").Value = '  -0.28%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.43'
$ws.Range("E20").Value = '  +1.41%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '193.56'
$ws.Range("E21").Value = '  -0.18%  '
$ws.Range("E22").Value = '  -0.17%  '
$ws.Range("E23").Value = '  +1.10%  '
$ws.Range("E24").Value = '  -0.25%  '
$ws.Range("E25").Value = '  -1.71%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '139.75'
$ws.Range("E26").Value = '  -1.21%  '
$ws.Range("E27").Value = '  -3.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.79'
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.46'
$ws.Range("E29").Value = '  -0.61%  '
$ws.Range("E30").Value = '  -1.43%  '
$ws.Range("E31").Value = '  -2.40%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.31'
$ws.Range("E32").Value = '  +0.30%  '
$ws.Range("E33").Value = '  +0.20%  '
$ws.Range("E34").Value = '  +0.66%  '
$ws.Range("E35").Value = '  -0.23%  '
$ws.Range("E36").Value = '  -1.37%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.56'
$ws.Range("E37").Value = '  -0.31%  '
$ws.Range("E38").Value = '  -2.03%  '
$ws.Range("D39").Value = '1.104.02'
$ws.Range("E39").Value = '  -2.57%  '
$ws.Range("E40").Value = '  -1.18%  '
$ws.Range("E41").Value = '  -0.24%  '
$ws.Range("E42").Value = '  +0.87%  '
$ws.Range("E43").Value = '  +0.81%  '
$ws.Range("E44").Value = '  -1.17%  '
$ws.Range("D45").Value = '1.762.26'
$ws.Range("E45").Value = '  -1.30%  '
$ws.Range("D46").Value = '0.0₆0108'
$ws.Range("E46").Value = '  -5.60%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '54.97'
$ws.Range("E47").Value = '  -1.36%  '
$ws.Range("B48").Value = 'SynthetixNetwork'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.39'
$ws.Range("E48").Value = '  +5.91%  '
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.418'
$ws.Range("E49").Value = '  -2.88%  '
$ws.Range("E50").Value = '  -1.22%  '
$ws.Range("E51").Value = '  -0.53%  '
